# Auto-generated Excel COM-interop script
# Applies the row-level numeric corrections described in the commit diff
# (scheduled runner refresh of computed Leve profit columns H:N)

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3995.0667
$ws.Range("I40").Value = 2705.5
$ws.Range("K40").Value = 2705.5
$ws.Range("M40").Value = -2530.5
# Row 62
$ws.Range("H62").Value = 3937.2727
$ws.Range("I62").Value = 3302.625
$ws.Range("J62").Value = 5629.6665
$ws.Range("K62").Value = 3302.625
$ws.Range("L62").Value = 5629.6665
$ws.Range("M62").Value = -2678.625
$ws.Range("N62").Value = -6877.6665
# Row 65
$ws.Range("H65").Value = 3937.2727
$ws.Range("I65").Value = 3302.625
$ws.Range("J65").Value = 5629.6665
$ws.Range("K65").Value = 16513.125
$ws.Range("L65").Value = 28148.3325
$ws.Range("M65").Value = -13393.125
$ws.Range("N65").Value = -34388.3325
# Row 76
$ws.Range("H76").Value = 4899.636
$ws.Range("I76").Value = 4737.375
$ws.Range("K76").Value = 4737.375
$ws.Range("M76").Value = -4422.375
# Row 79
$ws.Range("H79").Value = 4899.636
$ws.Range("I79").Value = 4737.375
$ws.Range("K79").Value = 4737.375
$ws.Range("M79").Value = -3645.375
# Row 82
$ws.Range("H82").Value = 12457.333
$ws.Range("I82").Value = 12457.333
$ws.Range("K82").Value = 37371.999
$ws.Range("M82").Value = -36965.999
# Row 85
$ws.Range("H85").Value = 12457.333
$ws.Range("I85").Value = 12457.333
$ws.Range("K85").Value = 37371.999
$ws.Range("M85").Value = -35967.999
# Row 98
$ws.Range("H98").Value = 3338.1904
$ws.Range("I98").Value = 3479.6843
$ws.Range("J98").Value = 1994
$ws.Range("K98").Value = 3479.6843
$ws.Range("L98").Value = 1994
$ws.Range("M98").Value = -1981.6843
$ws.Range("N98").Value = -4990
# Row 101
$ws.Range("H101").Value = 21709.777
$ws.Range("I101").Value = 1319.6
$ws.Range("J101").Value = 47197.5
$ws.Range("K101").Value = 3958.8
$ws.Range("L101").Value = 141592.5
$ws.Range("M101").Value = -2336.8
$ws.Range("N101").Value = -144836.5
# Row 106
$ws.Range("H106").Value = 1861.125
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 122
$ws.Range("H122").Value = 3338.1904
$ws.Range("I122").Value = 3479.6843
$ws.Range("J122").Value = 1994
$ws.Range("K122").Value = 10439.0529
$ws.Range("L122").Value = 5982
$ws.Range("M122").Value = -7989.052899999999
$ws.Range("N122").Value = -10882
# Row 137
$ws.Range("H137").Value = 2699.8
$ws.Range("I137").Value = 3963
$ws.Range("J137").Value = 2476.8823
$ws.Range("K137").Value = 11889
$ws.Range("L137").Value = 7430.646900000001
$ws.Range("M137").Value = -9339
$ws.Range("N137").Value = -12530.6469

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 2500
$ws.Range("J37").Value = 2500
$ws.Range("L37").Value = 2500
$ws.Range("N37").Value = -3046
# Row 61
$ws.Range("H61").Value = 1477.1428
$ws.Range("I61").Value = 1478.8
$ws.Range("J61").Value = 1473
$ws.Range("K61").Value = 1478.8
$ws.Range("L61").Value = 1473
$ws.Range("M61").Value = -1266.8
$ws.Range("N61").Value = -1897
# Row 63
$ws.Range("H63").Value = 23342.54
$ws.Range("I63").Value = 26732.182
$ws.Range("J63").Value = 4699.5
$ws.Range("K63").Value = 26732.182
$ws.Range("L63").Value = 4699.5
$ws.Range("M63").Value = -26046.182
$ws.Range("N63").Value = -6071.5
# Row 66
$ws.Range("H66").Value = 23342.54
$ws.Range("I66").Value = 26732.182
$ws.Range("J66").Value = 4699.5
$ws.Range("K66").Value = 133660.91
$ws.Range("L66").Value = 23497.5
$ws.Range("M66").Value = -130228.91
$ws.Range("N66").Value = -30361.5
# Row 88
$ws.Range("H88").Value = 1659.8462
$ws.Range("J88").Value = 2601.5
$ws.Range("L88").Value = 2601.5
$ws.Range("N88").Value = -3413.5
# Row 91
$ws.Range("H91").Value = 1659.8462
$ws.Range("J91").Value = 2601.5
$ws.Range("L91").Value = 2601.5
$ws.Range("N91").Value = -5409.5
# Row 97
$ws.Range("H97").Value = 5629.8667
$ws.Range("I97").Value = 6162.4165
$ws.Range("J97").Value = 3499.6667
$ws.Range("K97").Value = 6162.4165
$ws.Range("L97").Value = 3499.6667
$ws.Range("M97").Value = -5666.4165
$ws.Range("N97").Value = -4491.6667
# Row 102
$ws.Range("H102").Value = 1179.125
$ws.Range("I102").Value = 1188.8572
$ws.Range("K102").Value = 1188.8572
$ws.Range("M102").Value = 433.1428000000001
# Row 122
$ws.Range("H122").Value = 2472.476
$ws.Range("I122").Value = 2388.5334
$ws.Range("J122").Value = 2682.3333
$ws.Range("K122").Value = 7165.600199999999
$ws.Range("L122").Value = 8046.999899999999
$ws.Range("M122").Value = -4715.600199999999
$ws.Range("N122").Value = -12946.9999
# Row 136
$ws.Range("H136").Value = 1477.1428
$ws.Range("I136").Value = 1478.8
$ws.Range("J136").Value = 1473
$ws.Range("K136").Value = 4436.4
$ws.Range("L136").Value = 4419
$ws.Range("M136").Value = -1886.4
$ws.Range("N136").Value = -9519

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 7007
$ws.Range("I16").Value = 7007
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7007
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -6837
$ws.Range("N16").ClearContents()
# Row 20
$ws.Range("H20").Value = 7529.278
$ws.Range("I20").Value = 9240.77
$ws.Range("J20").Value = 3079.4
$ws.Range("K20").Value = 9240.77
$ws.Range("L20").Value = 3079.4
$ws.Range("M20").Value = -8993.77
$ws.Range("N20").Value = -3573.4
# Row 86
$ws.Range("H86").Value = 7862.9565
$ws.Range("I86").Value = 2604.7273
$ws.Range("K86").Value = 2604.7273
$ws.Range("M86").Value = -1481.7273
# Row 89
$ws.Range("H89").Value = 7862.9565
$ws.Range("I89").Value = 2604.7273
$ws.Range("K89").Value = 13023.6365
$ws.Range("M89").Value = -7407.636500000001
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 715.2
$ws.Range("I10").Value = 761.5
$ws.Range("J10").Value = 530
$ws.Range("K10").Value = 761.5
$ws.Range("L10").Value = 530
$ws.Range("M10").Value = -622.5
$ws.Range("N10").Value = -808
# Row 62
$ws.Range("H62").Value = 36766.39
$ws.Range("I62").Value = 3673.0667
$ws.Range("K62").Value = 3673.0667
$ws.Range("M62").Value = -3049.0667
# Row 65
$ws.Range("H65").Value = 36766.39
$ws.Range("I65").Value = 3673.0667
$ws.Range("K65").Value = 18365.3335
$ws.Range("M65").Value = -15245.3335
# Row 86
$ws.Range("H86").Value = 2799.6155
$ws.Range("I86").Value = 2575.5
$ws.Range("J86").Value = 3158.2
$ws.Range("K86").Value = 2575.5
$ws.Range("L86").Value = 3158.2
$ws.Range("M86").Value = -1452.5
$ws.Range("N86").Value = -5404.2
# Row 89
$ws.Range("H89").Value = 2799.6155
$ws.Range("I89").Value = 2575.5
$ws.Range("J89").Value = 3158.2
$ws.Range("K89").Value = 12877.5
$ws.Range("L89").Value = 15791
$ws.Range("M89").Value = -7261.5
$ws.Range("N89").Value = -27023

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 655.2
$ws.Range("I103").Value = 294
$ws.Range("K103").Value = 882
$ws.Range("M103").Value = -3
# Row 131
$ws.Range("H131").Value = 272929.75
$ws.Range("I131").Value = 913.3333
$ws.Range("J131").Value = 299685.47
$ws.Range("K131").Value = 2739.9999
$ws.Range("L131").Value = 899056.4099999999
$ws.Range("M131").Value = 2300.0001
$ws.Range("N131").Value = -909136.4099999999
# Row 132
$ws.Range("H132").Value = 9155.299999999999
$ws.Range("I132").Value = 3238.8
$ws.Range("J132").Value = 15071.8
$ws.Range("K132").Value = 29149.2
$ws.Range("L132").Value = 135646.2
$ws.Range("M132").Value = -26619.2
$ws.Range("N132").Value = -140706.2

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6906.225
$ws.Range("I70").Value = 4999.6665
$ws.Range("J70").Value = 7060.811
$ws.Range("K70").Value = 4999.6665
$ws.Range("L70").Value = 7060.811
$ws.Range("M70").Value = -4729.6665
$ws.Range("N70").Value = -7600.811
# Row 73
$ws.Range("H73").Value = 6906.225
$ws.Range("I73").Value = 4999.6665
$ws.Range("J73").Value = 7060.811
$ws.Range("K73").Value = 4999.6665
$ws.Range("L73").Value = 7060.811
$ws.Range("M73").Value = -4063.6665
$ws.Range("N73").Value = -8932.811
# Row 80
$ws.Range("H80").Value = 4138.4287
$ws.Range("I80").Value = 2401.1667
$ws.Range("J80").Value = 4833.3335
$ws.Range("K80").Value = 2401.1667
$ws.Range("L80").Value = 4833.3335
$ws.Range("M80").Value = -1403.1667
$ws.Range("N80").Value = -6829.3335
# Row 83
$ws.Range("H83").Value = 4138.4287
$ws.Range("I83").Value = 2401.1667
$ws.Range("J83").Value = 4833.3335
$ws.Range("K83").Value = 12005.8335
$ws.Range("L83").Value = 24166.6675
$ws.Range("M83").Value = -7013.833500000001
$ws.Range("N83").Value = -34150.6675
# Row 132
$ws.Range("H132").Value = 1513.8462
$ws.Range("I132").Value = 1399.826
$ws.Range("J132").Value = 2388
$ws.Range("K132").Value = 4199.478
$ws.Range("L132").Value = 7164
$ws.Range("M132").Value = -1669.478
$ws.Range("N132").Value = -12224

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2078.6667
$ws.Range("J68").Value = 2166.3333
$ws.Range("L68").Value = 2166.3333
$ws.Range("N68").Value = -3664.3333
# Row 71
$ws.Range("H71").Value = 2078.6667
$ws.Range("J71").Value = 2166.3333
$ws.Range("L71").Value = 10831.6665
$ws.Range("N71").Value = -18319.6665
# Row 76
$ws.Range("H76").Value = 9990
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 9990
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 9990
$ws.Range("N76").Value = -10666
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 9990
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 9990
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 9990
$ws.Range("N79").Value = -12330
$ws.Range("M79").ClearContents()
# Row 132
$ws.Range("H132").Value = 4507.7
$ws.Range("I132").Value = 2787.3333
$ws.Range("K132").Value = 8361.999899999999
$ws.Range("M132").Value = -5831.999899999999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 141
$ws.Range("H141").Value = 99500
$ws.Range("J141").Value = 99500
$ws.Range("L141").Value = 99500
$ws.Range("N141").Value = -109860

